# Generate Report for Handback
#
# This script mutates the localization-status workbook the way the
# "Generate Report for Handback" job would: it flips the status text from
# "Ready for handoff" to "Handed back: in sync with en-US", records the
# handback timestamps for the zh-cn and de-de jobs, fills in the
# "Latest Target File" (as a hyperlink) / "Latest Handback File" columns,
# and widens a few columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Flip every "Ready for handoff" status cell to the handed-back text.
# ---------------------------------------------------------------------
foreach ($cellRef in @("E2", "F2", "E3", "F3")) {
    $cell = $wsOverview.Range($cellRef)
    if ($cell.Value() -eq $oldStatus) {
        $cell.Value = $newStatus
    }
}

foreach ($ws in @($wsZhCn, $wsDeDe)) {
    foreach ($cellRef in @("C2", "C3")) {
        $cell = $ws.Range($cellRef)
        if ($cell.Value() -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# ---------------------------------------------------------------------
# 2. zh-cn handback: target file hyperlink, handback xliff name, and the
#    handback timestamp.
# ---------------------------------------------------------------------
$zhCnTargetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/572db8f938289597fb22ba28d2b976003814b23f/e2e/e1ef11d8-2e75-4ed4-b1d4-b8e55d6bb396.md"
$zhCnTargetDisplay = "e1ef11d8-2e75-4ed4-b1d4-b8e55d6bb396.md"
$zhCnHandbackFile = "e1ef11d8-2e75-4ed4-b1d4-b8e55d6bb396.49819bff66962a9e3321e25b8aa4e76d9a5ccadb.zh-cn.xlf"
$zhCnHandbackDate = "2016-08-28 05:06:09"

foreach ($rowNum in @(2, 3)) {
    $targetCell = $wsZhCn.Range("I" + $rowNum)
    $wsZhCn.Hyperlinks.Add($targetCell, $zhCnTargetUrl, "", "", $zhCnTargetDisplay)
    $targetCell.Font.Underline = $true
    $targetCell.Font.Color = 15570276

    $wsZhCn.Range("J" + $rowNum).Value = $zhCnHandbackFile
    $wsZhCn.Range("K" + $rowNum).Value = $zhCnHandbackDate
}

# ---------------------------------------------------------------------
# 3. de-de handback: target file hyperlink, handback xliff name, and the
#    (later) handback timestamp.
# ---------------------------------------------------------------------
$deDeTargetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/572db8f938289597fb22ba28d2b976003814b23f/e2e/e1ef11d8-2e75-4ed4-b1d4-b8e55d6bb396.md"
$deDeTargetDisplay = "e1ef11d8-2e75-4ed4-b1d4-b8e55d6bb396.md"
$deDeHandbackFile = "e1ef11d8-2e75-4ed4-b1d4-b8e55d6bb396.49819bff66962a9e3321e25b8aa4e76d9a5ccadb.de-de.xlf"
$deDeHandbackDate = "2016-08-28 05:06:16"

foreach ($rowNum in @(2, 3)) {
    $targetCell = $wsDeDe.Range("I" + $rowNum)
    $wsDeDe.Hyperlinks.Add($targetCell, $deDeTargetUrl, "", "", $deDeTargetDisplay)
    $targetCell.Font.Underline = $true
    $targetCell.Font.Color = 15570276

    $wsDeDe.Range("J" + $rowNum).Value = $deDeHandbackFile
    $wsDeDe.Range("K" + $rowNum).Value = $deDeHandbackDate
}

# ---------------------------------------------------------------------
# 4. Widen the columns that now hold the longer status / file-name text.
#    (ColumnWidth is quantized by Excel to 1/7-character pixel steps, so
#    these inputs are chosen to land on the nearest achievable width.)
# ---------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.16
$wsOverview.Columns.Item(6).ColumnWidth = 29.16

foreach ($ws in @($wsZhCn, $wsDeDe)) {
    $ws.Columns.Item(3).ColumnWidth = 29.16
    $ws.Columns.Item(9).ColumnWidth = 39.16
    $ws.Columns.Item(10).ColumnWidth = 39.16
}
